$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing rows 2-5: sending cluster changes from FAPs to ECs, and values are refreshed
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lrp6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06937033333333333
$ws.Range("H2").Value = 0.208111
$ws.Range("I2").Value = 0.01708561286819356
$ws.Range("J2").Value = 0.01708561286819356
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.779612
$ws.Range("N2").Value = 32.338836
$ws.Range("O2").Value = 0.1321092878737708
$ws.Range("P2").Value = 0.1321092878737708
$ws.Range("Q2").Value = 0.7477852776439999
$ws.Range("R2").Value = 6.730067498796
$ws.Range("S2").Value = 0.002257168148903985
$ws.Range("T2").Value = 0.002257168148903985

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lrp6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06937033333333333
$ws.Range("H3").Value = 0.208111
$ws.Range("I3").Value = 0.01708561286819356
$ws.Range("J3").Value = 0.01708561286819356
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 38.54369466666667
$ws.Range("N3").Value = 115.631084
$ws.Range("O3").Value = 0.4723713668393066
$ws.Range("P3").Value = 0.4723713668393065
$ws.Range("Q3").Value = 2.673788946924889
$ws.Range("R3").Value = 24.064100522324
$ws.Range("S3").Value = 0.008070754303835834
$ws.Range("T3").Value = 0.008070754303835836

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lrp6"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06937033333333333
$ws.Range("H4").Value = 0.208111
$ws.Range("I4").Value = 0.01708561286819356
$ws.Range("J4").Value = 0.01708561286819356
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 12.62567333333333
$ws.Range("N4").Value = 37.87702
$ws.Range("O4").Value = 0.1547336502458089
$ws.Range("P4").Value = 0.1547336502458089
$ws.Range("Q4").Value = 0.875847167691111
$ws.Range("R4").Value = 7.88262450922
$ws.Range("S4").Value = 0.002643719245782354
$ws.Range("T4").Value = 0.002643719245782354

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lrp6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06937033333333333
$ws.Range("H5").Value = 0.208111
$ws.Range("I5").Value = 0.01708561286819356
$ws.Range("J5").Value = 0.01708561286819356
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.64719066666667
$ws.Range("N5").Value = 58.94157200000001
$ws.Range("O5").Value = 0.2407856950411137
$ws.Range("P5").Value = 0.2407856950411137
$ws.Range("Q5").Value = 1.362932165610222
$ws.Range("R5").Value = 12.266389490492
$ws.Range("S5").Value = 0.004113971169671382
$ws.Range("T5").Value = 0.004113971169671382

# Add new rows 6-9 with sending cluster FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Lrp6"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.990790333333333
$ws.Range("H6").Value = 11.972371
$ws.Range("I6").Value = 0.9829143871318063
$ws.Range("J6").Value = 0.9829143871318063
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.779612
$ws.Range("N6").Value = 32.338836
$ws.Range("O6").Value = 0.1321092878737708
$ws.Range("P6").Value = 0.1321092878737708
$ws.Range("Q6").Value = 43.019171366684
$ws.Range("R6").Value = 387.172542300156
$ws.Range("S6").Value = 0.1298521197248668
$ws.Range("T6").Value = 0.1298521197248668

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Lrp6"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.990790333333333
$ws.Range("H7").Value = 11.972371
$ws.Range("I7").Value = 0.9829143871318063
$ws.Range("J7").Value = 0.9829143871318063
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 38.54369466666667
$ws.Range("N7").Value = 115.631084
$ws.Range("O7").Value = 0.4723713668393066
$ws.Range("P7").Value = 0.4723713668393065
$ws.Range("Q7").Value = 153.8198040866849
$ws.Range("R7").Value = 1384.378236780164
$ws.Range("S7").Value = 0.4643006125354707
$ws.Range("T7").Value = 0.4643006125354707

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Rspo3"
$ws.Range("C8").Value = "Lrp6"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.990790333333333
$ws.Range("H8").Value = 11.972371
$ws.Range("I8").Value = 0.9829143871318063
$ws.Range("J8").Value = 0.9829143871318063
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.62567333333333
$ws.Range("N8").Value = 37.87702
$ws.Range("O8").Value = 0.1547336502458089
$ws.Range("P8").Value = 0.1547336502458089
$ws.Range("Q8").Value = 50.38641509049111
$ws.Range("R8").Value = 453.47773581442
$ws.Range("S8").Value = 0.1520899310000265
$ws.Range("T8").Value = 0.1520899310000265

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Rspo3"
$ws.Range("C9").Value = "Lrp6"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.990790333333333
$ws.Range("H9").Value = 11.972371
$ws.Range("I9").Value = 0.9829143871318063
$ws.Range("J9").Value = 0.9829143871318063
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.64719066666667
$ws.Range("N9").Value = 58.94157200000001
$ws.Range("O9").Value = 0.2407856950411137
$ws.Range("P9").Value = 0.2407856950411137
$ws.Range("Q9").Value = 78.40781858969024
$ws.Range("R9").Value = 705.670367307212
$ws.Range("S9").Value = 0.2366717238714423
$ws.Range("T9").Value = 0.2366717238714423
